# Apply case-1 data edit: refresh column A sample values + widen column A slightly
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A width: stored OOXML width goes from 15.42578125 -> 15.7109375.
# 15.7109375 is exactly what Excel stores when ColumnWidth is set to 14.8
# (closest setting reachable through this host's ColumnWidth quantization).
$ws.Columns.Item(1).ColumnWidth = 14.8

# Refresh the 33 sampled values in column A
$ws.Range("A1").Value = 0.39698209129538498
$ws.Range("A2").Value = -0.0099999995871904446
$ws.Range("A3").Value = -0.034319509131092829
$ws.Range("A4").Value = -0.011999999886414869
$ws.Range("A5").Value = -0.0059999995859918442
$ws.Range("A6").Value = -0.0059999995743673651
$ws.Range("A7").Value = -0.019999999496041809
$ws.Range("A8").Value = -0.01999999949413489
$ws.Range("A9").Value = 0.034594747688743155
$ws.Range("A10").Value = -0.0059999995658586158
$ws.Range("A11").Value = -0.0044999995736070275
$ws.Range("A12").Value = -0.0059999995642496806
$ws.Range("A13").Value = -0.0059999995598918332
$ws.Range("A14").Value = -0.01199999952550268
$ws.Range("A15").Value = -0.0059999995580160004
$ws.Range("A16").Value = 0.0080604992788781971
$ws.Range("A17").Value = -0.0059999995557626917
$ws.Range("A18").Value = -0.0089999995386929044
$ws.Range("A19").Value = -0.0089999995914507025
$ws.Range("A20").Value = -0.066646365703947552
$ws.Range("A21").Value = -0.008999999578152007
$ws.Range("A22").Value = -0.0089999995777114705
$ws.Range("A23").Value = -0.008999999575443951
$ws.Range("A24").Value = -0.041999999386234776
$ws.Range("A25").Value = -0.041999999382934305
$ws.Range("A26").Value = -0.0059999995730883882
$ws.Range("A27").Value = -0.0059999995712214371
$ws.Range("A28").Value = -0.0059999995651001115
$ws.Range("A29").Value = -0.011999999527146699
$ws.Range("A30").Value = -0.019999999480686981
$ws.Range("A31").Value = 0.00036892843383640184
$ws.Range("A32").Value = -0.02099999947188369
$ws.Range("A33").Value = -0.0059999995546773377
